$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the "R10" rule (row 8, col E) to reflect the git update
$ws.Range("E8").Value = "GIT UPDATE"

# Select cell E8, mirroring the active selection recorded in the saved file
$ws.Activate()
$ws.Range("E8").Select()
